$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing nothing (row 4 currently empty),
# then copy the existing row 3 (old weekly entry) down into row 4 before
# overwriting row 3 with the new weekly entry.
$ws.Rows.Item(4).Insert()

# Copy original row 3 content into the newly inserted row 4
$ws.Range("A3:R3").Copy()
$ws.Range("A4:R4").PasteSpecial()

# Update row 3 with the new weekly values
$ws.Range("D3").Value = 44775
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 8000
$ws.Range("M3").Value = 8000
$ws.Range("P3").Value = 800
